$d = $word.ActiveDocument

# Locate the two target paragraphs by their current text instead of a fixed
# index, so the script is resilient to the exact paragraph numbering.
$p1 = $null
$p3 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $para = $d.Paragraphs($i)
  $t = $para.Range.Text
  if ($p1 -eq $null -and $t -like "*dificuldade de gest*") {
    $p1 = $para
  }
  if ($p3 -eq $null -and $t -like "*Criar um fluxo*") {
    $p3 = $para
  }
}
if ($p1 -eq $null) { $p1 = $d.Paragraphs(3) }
if ($p3 -eq $null) { $p3 = $d.Paragraphs(6) }

# ---------------------------------------------------------------------------
# 1) First body paragraph: replace the NormalWeb-styled single run with
#    explicit run formatting split across three runs (plus a proofErr pair
#    around the stray "a") and move rPr from a removed pPr into the pPr.
# ---------------------------------------------------------------------------
$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="36F3FAC6" w14:textId="77777777" w:rsidR="008A1ADF" w:rsidRDefault="008A1ADF" w:rsidP="008A1ADF">
            <w:pPr>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:color w:val="202124"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:color w:val="202124"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t xml:space="preserve">O problema da dificuldade de gestão da escola, afeta ao gestor e secretaria, devido </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:color w:val="202124"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t>a</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:color w:val="202124"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
                <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              </w:rPr>
              <w:t xml:space="preserve"> demora para realizar processos, falhas no controle das turmas e dos cronogramas e dificuldade na tomada de decisões.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Empty paragraph right after: drop the stray "_GoBack" bookmark pair.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 3) First bullet item: split "Criar um fluxo para cadastramento de alunos;"
#    into multiple runs and re-insert a "_GoBack" bookmark around "fluxo ".
# ---------------------------------------------------------------------------
$xml3 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="77DA2A5F" w14:textId="77777777" w:rsidR="00D25E7E" w:rsidRDefault="00D25E7E" w:rsidP="00D25E7E">
            <w:pPr>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="2"/>
              </w:numPr>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Criar um </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>fluxo</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>para cadastramento de alunos;</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$p3.Range.InsertXML($xml3)

Write-Output "done"
